# Weekly update: a new price record is inserted at row 72 for this
# "Hortaliza, Femacal de La Calera - Haba" subset, pushing the previously
# existing rows 72-79 down to 73-80 (dimension grows from R79 to R80).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 72..79 down to 73..80, leaving a fresh blank row 72.
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new weekly record
# (Provincia de Quillota).
$ws.Range("A72").Value = 3
$ws.Range("B72").Value = "Femacal de La Calera"
$ws.Range("C72").Value = "Coquimbo"
$ws.Range("D72").Value = 44491
$ws.Range("E72").Value = 5
$ws.Range("F72").Value = 100112026
$ws.Range("G72").Value = "Haba"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 110
$ws.Range("K72").Value = 8000
$ws.Range("L72").Value = 8500
$ws.Range("M72").Value = 8227
$ws.Range("N72").Value = "$/saco 25 kilos"
$ws.Range("O72").Value = "Provincia de Quillota"
$ws.Range("P72").Value = 329
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"

# The diff also shows N78 ($/saco -> $/malla) differing from a pure
# row-shift of the old data (every other column for that row matches a
# straight shift). Apply that extra correction explicitly.
$ws.Range("N78").Value = "$/malla 25 kilos"
